$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A3:D113").Clear()
$ws.Range("A3").Value = 44235
$ws.Range("B3").Value = "Wages Expense"
$ws.Range("E2").Value = 934525
